# Balance sheet yearly update: roll the 5-year window forward one year
# (drop FY1396, add FY1401) and refresh the "as of" publish dates.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 8: period headers (12-month period ended ...) ---
$ws.Cells.Item(8,4).Value = "12 ماهه منتهی به 1397/12"
$ws.Cells.Item(8,5).Value = "12 ماهه منتهی به 1398/12"
$ws.Cells.Item(8,6).Value = "12 ماهه منتهی به 1399/12"
$ws.Cells.Item(8,7).Value = "12 ماهه منتهی به 1400/12"
$ws.Cells.Item(8,8).Value = "12 ماهه منتهی به 1401/12"

# --- Row 9: publish dates ---
$ws.Cells.Item(9,4).Value = "1399-03-24 (10)"
$ws.Cells.Item(9,5).Value = "1400-04-02 (9)"
$ws.Cells.Item(9,6).Value = "1401-04-08 (10)"
$ws.Cells.Item(9,7).Value = "1402-01-29 (8)"
$ws.Cells.Item(9,8).Value = "1402-01-29"

# --- Rows 12-58: balance-sheet figures, shifted one column left (FY1396 dropped) with FY1401 figures filling column H ---

# Row 12
$ws.Cells.Item(12,4).Value = 3874497
$ws.Cells.Item(12,5).Value = 6479365
$ws.Cells.Item(12,6).Value = 108790068
$ws.Cells.Item(12,7).Value = 99912743
$ws.Cells.Item(12,8).Value = 75240206

# Row 13
$ws.Cells.Item(13,4).Value = 0
$ws.Cells.Item(13,5).Value = 0
$ws.Cells.Item(13,6).Value = 2531507
$ws.Cells.Item(13,7).Value = 194520
$ws.Cells.Item(13,8).Value = 139894

# Row 14
$ws.Cells.Item(14,4).Value = 9970073
$ws.Cells.Item(14,5).Value = 13102323
$ws.Cells.Item(14,6).Value = 50569013
$ws.Cells.Item(14,7).Value = 107696139
$ws.Cells.Item(14,8).Value = 113927165

# Row 15
$ws.Cells.Item(15,4).Value = 14279916
$ws.Cells.Item(15,5).Value = 24685807
$ws.Cells.Item(15,6).Value = 31594275
$ws.Cells.Item(15,7).Value = 50057714
$ws.Cells.Item(15,8).Value = 78756902

# Row 16
$ws.Cells.Item(16,4).Value = 9190087
$ws.Cells.Item(16,5).Value = 13222327
$ws.Cells.Item(16,6).Value = 13988958
$ws.Cells.Item(16,7).Value = 26540986
$ws.Cells.Item(16,8).Value = 32153400

# Row 17
$ws.Cells.Item(17,4).Value = 0
$ws.Cells.Item(17,5).Value = 0
$ws.Cells.Item(17,6).Value = 0
$ws.Cells.Item(17,7).Value = 0
$ws.Cells.Item(17,8).Value = 0

# Row 18
$ws.Cells.Item(18,4).Value = 37314573
$ws.Cells.Item(18,5).Value = 57489822
$ws.Cells.Item(18,6).Value = 207473821
$ws.Cells.Item(18,7).Value = 284402102
$ws.Cells.Item(18,8).Value = 300217567

# Row 19
$ws.Cells.Item(19,4).Value = 0
$ws.Cells.Item(19,5).Value = 0
$ws.Cells.Item(19,6).Value = 0
$ws.Cells.Item(19,7).Value = 6490291
$ws.Cells.Item(19,8).Value = 15916608

# Row 20
$ws.Cells.Item(20,4).Value = 18309420
$ws.Cells.Item(20,5).Value = 29686063
$ws.Cells.Item(20,6).Value = 51931045
$ws.Cells.Item(20,7).Value = 125391256
$ws.Cells.Item(20,8).Value = 173628112

# Row 21
$ws.Cells.Item(21,4).Value = 0
$ws.Cells.Item(21,5).Value = 0
$ws.Cells.Item(21,6).Value = 0
$ws.Cells.Item(21,7).Value = 0
$ws.Cells.Item(21,8).Value = 0

# Row 22
$ws.Cells.Item(22,4).Value = 43215135
$ws.Cells.Item(22,5).Value = 42773519
$ws.Cells.Item(22,6).Value = 46104801
$ws.Cells.Item(22,7).Value = 73603677
$ws.Cells.Item(22,8).Value = 83844833

# Row 23
$ws.Cells.Item(23,4).Value = 42893
$ws.Cells.Item(23,5).Value = 32131
$ws.Cells.Item(23,6).Value = 24098
$ws.Cells.Item(23,7).Value = 39605
$ws.Cells.Item(23,8).Value = 266083

# Row 24
$ws.Cells.Item(24,4).Value = "-"
$ws.Cells.Item(24,5).Value = "-"
$ws.Cells.Item(24,6).Value = "-"
$ws.Cells.Item(24,7).Value = "-"
$ws.Cells.Item(24,8).Value = "-"

# Row 25
$ws.Cells.Item(25,4).Value = 0
$ws.Cells.Item(25,5).Value = 0
$ws.Cells.Item(25,6).Value = 0
$ws.Cells.Item(25,7).Value = 0
$ws.Cells.Item(25,8).Value = 0

# Row 26
$ws.Cells.Item(26,4).Value = 61567448
$ws.Cells.Item(26,5).Value = 72491713
$ws.Cells.Item(26,6).Value = 98059944
$ws.Cells.Item(26,7).Value = 205524829
$ws.Cells.Item(26,8).Value = 273655636

# Row 27
$ws.Cells.Item(27,4).Value = 98882021
$ws.Cells.Item(27,5).Value = 129981535
$ws.Cells.Item(27,6).Value = 305533765
$ws.Cells.Item(27,7).Value = 489926931
$ws.Cells.Item(27,8).Value = 573873203

# Row 29
$ws.Cells.Item(29,4).Value = 18431386
$ws.Cells.Item(29,5).Value = 13818118
$ws.Cells.Item(29,6).Value = 36651410
$ws.Cells.Item(29,7).Value = 49655148
$ws.Cells.Item(29,8).Value = 66718107

# Row 30
$ws.Cells.Item(30,4).Value = "-"
$ws.Cells.Item(30,5).Value = "-"
$ws.Cells.Item(30,6).Value = "-"
$ws.Cells.Item(30,7).Value = "-"
$ws.Cells.Item(30,8).Value = "-"

# Row 31
$ws.Cells.Item(31,4).Value = 1362125
$ws.Cells.Item(31,5).Value = 421186
$ws.Cells.Item(31,6).Value = 2142550
$ws.Cells.Item(31,7).Value = 5991138
$ws.Cells.Item(31,8).Value = 6864183

# Row 32
$ws.Cells.Item(32,4).Value = 0
$ws.Cells.Item(32,5).Value = 9300772
$ws.Cells.Item(32,6).Value = 14692100
$ws.Cells.Item(32,7).Value = 38418943
$ws.Cells.Item(32,8).Value = 17728384

# Row 33
$ws.Cells.Item(33,4).Value = 2188728
$ws.Cells.Item(33,5).Value = 437369
$ws.Cells.Item(33,6).Value = 6811586
$ws.Cells.Item(33,7).Value = 702542
$ws.Cells.Item(33,8).Value = 413735

# Row 34
$ws.Cells.Item(34,4).Value = 3754613
$ws.Cells.Item(34,5).Value = 1186016
$ws.Cells.Item(34,6).Value = 0
$ws.Cells.Item(34,7).Value = 0
$ws.Cells.Item(34,8).Value = 0

# Row 35
$ws.Cells.Item(35,4).Value = 0
$ws.Cells.Item(35,5).Value = 0
$ws.Cells.Item(35,6).Value = 0
$ws.Cells.Item(35,7).Value = 0
$ws.Cells.Item(35,8).Value = 0

# Row 36
$ws.Cells.Item(36,4).Value = 0
$ws.Cells.Item(36,5).Value = 0
$ws.Cells.Item(36,6).Value = 0
$ws.Cells.Item(36,7).Value = 0
$ws.Cells.Item(36,8).Value = 0

# Row 37
$ws.Cells.Item(37,4).Value = 25736852
$ws.Cells.Item(37,5).Value = 25163461
$ws.Cells.Item(37,6).Value = 60297646
$ws.Cells.Item(37,7).Value = 94767771
$ws.Cells.Item(37,8).Value = 91724409

# Row 38
$ws.Cells.Item(38,4).Value = 152474
$ws.Cells.Item(38,5).Value = 0
$ws.Cells.Item(38,6).Value = 0
$ws.Cells.Item(38,7).Value = 0
$ws.Cells.Item(38,8).Value = 0

# Row 39
$ws.Cells.Item(39,4).Value = "-"
$ws.Cells.Item(39,5).Value = "-"
$ws.Cells.Item(39,6).Value = "-"
$ws.Cells.Item(39,7).Value = "-"
$ws.Cells.Item(39,8).Value = "-"

# Row 40
$ws.Cells.Item(40,4).Value = 0
$ws.Cells.Item(40,5).Value = 0
$ws.Cells.Item(40,6).Value = 0
$ws.Cells.Item(40,7).Value = 0
$ws.Cells.Item(40,8).Value = 0

# Row 41
$ws.Cells.Item(41,4).Value = 147500
$ws.Cells.Item(41,5).Value = 178244
$ws.Cells.Item(41,6).Value = 194397
$ws.Cells.Item(41,7).Value = 314902
$ws.Cells.Item(41,8).Value = 486313

# Row 42
$ws.Cells.Item(42,4).Value = 299974
$ws.Cells.Item(42,5).Value = 178244
$ws.Cells.Item(42,6).Value = 194397
$ws.Cells.Item(42,7).Value = 314902
$ws.Cells.Item(42,8).Value = 486313

# Row 43
$ws.Cells.Item(43,4).Value = 26036826
$ws.Cells.Item(43,5).Value = 25341705
$ws.Cells.Item(43,6).Value = 60492043
$ws.Cells.Item(43,7).Value = 95082673
$ws.Cells.Item(43,8).Value = 92210722

# Row 45
$ws.Cells.Item(45,4).Value = 33500000
$ws.Cells.Item(45,5).Value = 44000000
$ws.Cells.Item(45,6).Value = 55500000
$ws.Cells.Item(45,7).Value = 119000000
$ws.Cells.Item(45,8).Value = 192500000

# Row 46
$ws.Cells.Item(46,4).Value = 0
$ws.Cells.Item(46,5).Value = 0
$ws.Cells.Item(46,6).Value = 0
$ws.Cells.Item(46,7).Value = 0
$ws.Cells.Item(46,8).Value = 0

# Row 47
$ws.Cells.Item(47,4).Value = 7053315
$ws.Cells.Item(47,5).Value = 10595020
$ws.Cells.Item(47,6).Value = 0
$ws.Cells.Item(47,7).Value = 0
$ws.Cells.Item(47,8).Value = 71613156

# Row 48
$ws.Cells.Item(48,4).Value = 0
$ws.Cells.Item(48,5).Value = 0
$ws.Cells.Item(48,6).Value = -698575
$ws.Cells.Item(48,7).Value = -2203228
$ws.Cells.Item(48,8).Value = -2567191

# Row 49
$ws.Cells.Item(49,4).Value = 0
$ws.Cells.Item(49,5).Value = 0
$ws.Cells.Item(49,6).Value = 62707
$ws.Cells.Item(49,7).Value = 108002
$ws.Cells.Item(49,8).Value = 97945

# Row 50
$ws.Cells.Item(50,4).Value = 3350000
$ws.Cells.Item(50,5).Value = 4400000
$ws.Cells.Item(50,6).Value = 5550000
$ws.Cells.Item(50,7).Value = 11900000
$ws.Cells.Item(50,8).Value = 19250000

# Row 51
$ws.Cells.Item(51,4).Value = 0
$ws.Cells.Item(51,5).Value = 0
$ws.Cells.Item(51,6).Value = 0
$ws.Cells.Item(51,7).Value = 0
$ws.Cells.Item(51,8).Value = 0

# Row 52
$ws.Cells.Item(52,4).Value = "-"
$ws.Cells.Item(52,5).Value = "-"
$ws.Cells.Item(52,6).Value = "-"
$ws.Cells.Item(52,7).Value = "-"
$ws.Cells.Item(52,8).Value = "-"

# Row 53
$ws.Cells.Item(53,4).Value = 0
$ws.Cells.Item(53,5).Value = 0
$ws.Cells.Item(53,6).Value = 0
$ws.Cells.Item(53,7).Value = 0
$ws.Cells.Item(53,8).Value = 0

# Row 54
$ws.Cells.Item(54,4).Value = "-"
$ws.Cells.Item(54,5).Value = "-"
$ws.Cells.Item(54,6).Value = "-"
$ws.Cells.Item(54,7).Value = "-"
$ws.Cells.Item(54,8).Value = "-"

# Row 55
$ws.Cells.Item(55,4).Value = 0
$ws.Cells.Item(55,5).Value = 0
$ws.Cells.Item(55,6).Value = 0
$ws.Cells.Item(55,7).Value = 0
$ws.Cells.Item(55,8).Value = 0

# Row 56
$ws.Cells.Item(56,4).Value = 28941880
$ws.Cells.Item(56,5).Value = 45644810
$ws.Cells.Item(56,6).Value = 184627590
$ws.Cells.Item(56,7).Value = 266039484
$ws.Cells.Item(56,8).Value = 200768571

# Row 57
$ws.Cells.Item(57,4).Value = 72845195
$ws.Cells.Item(57,5).Value = 104639830
$ws.Cells.Item(57,6).Value = 245041722
$ws.Cells.Item(57,7).Value = 394844258
$ws.Cells.Item(57,8).Value = 481662481

# Row 58
$ws.Cells.Item(58,4).Value = 98882021
$ws.Cells.Item(58,5).Value = 129981535
$ws.Cells.Item(58,6).Value = 305533765
$ws.Cells.Item(58,7).Value = 489926931
$ws.Cells.Item(58,8).Value = 573873203

